# Applies scheduled market-data refresh updates to the Lamia_Profits workbook.
# Updates currentAveragePrice* / Leve price / profit columns (H:N) across all
# job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the latest data pull.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 109.125
$ws.Range("I2").Value = 85.5
$ws.Range("J2").Value = 180
$ws.Range("K2").Value = 85.5
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = 27.5
$ws.Range("N2").Value = -406
$ws.Range("H15").Value = 1444.8529
$ws.Range("I15").Value = 1444.8529
$ws.Range("K15").Value = 4334.5587
$ws.Range("M15").Value = -4165.5587
$ws.Range("H19").Value = 639.3333
$ws.Range("I19").Value = 639.75
$ws.Range("J19").Value = 639
$ws.Range("K19").Value = 639.75
$ws.Range("L19").Value = 639
$ws.Range("M19").Value = -464.75
$ws.Range("N19").Value = -989
$ws.Range("H40").Value = 5259.8887
$ws.Range("I40").Value = 3675
$ws.Range("J40").Value = 5712.7144
$ws.Range("K40").Value = 3675
$ws.Range("L40").Value = 5712.7144
$ws.Range("M40").Value = -3500
$ws.Range("N40").Value = -6062.7144
$ws.Range("H41").Value = 1454.7
$ws.Range("I41").Value = 1781.8125
$ws.Range("K41").Value = 1781.8125
$ws.Range("M41").Value = -1341.8125
$ws.Range("H55").Value = 546.875
$ws.Range("I55").Value = 348.75
$ws.Range("J55").Value = 943.125
$ws.Range("K55").Value = 348.75
$ws.Range("L55").Value = 943.125
$ws.Range("M55").Value = -134.75
$ws.Range("N55").Value = -1371.125
$ws.Range("H74").Value = 11192.375
$ws.Range("I74").Value = 11493.667
$ws.Range("J74").Value = 11011.6
$ws.Range("K74").Value = 11493.667
$ws.Range("L74").Value = 11011.6
$ws.Range("M74").Value = -10557.667
$ws.Range("N74").Value = -12883.6
$ws.Range("H77").Value = 11192.375
$ws.Range("I77").Value = 11493.667
$ws.Range("J77").Value = 11011.6
$ws.Range("K77").Value = 57468.335
$ws.Range("L77").Value = 55058
$ws.Range("M77").Value = -52788.335
$ws.Range("N77").Value = -64418
$ws.Range("H80").Value = 2023.2106
$ws.Range("I80").Value = 572.9231
$ws.Range("K80").Value = 1718.7693
$ws.Range("M80").Value = -720.7692999999999
$ws.Range("H83").Value = 2023.2106
$ws.Range("I83").Value = 572.9231
$ws.Range("K83").Value = 5156.3079
$ws.Range("M83").Value = -164.3078999999998
$ws.Range("H86").Value = 2095.9333
$ws.Range("I86").Value = 1873
$ws.Range("J86").Value = 2244.5557
$ws.Range("K86").Value = 1873
$ws.Range("L86").Value = 2244.5557
$ws.Range("M86").Value = -750
$ws.Range("N86").Value = -4490.5557
$ws.Range("H89").Value = 2095.9333
$ws.Range("I89").Value = 1873
$ws.Range("J89").Value = 2244.5557
$ws.Range("K89").Value = 9365
$ws.Range("L89").Value = 11222.7785
$ws.Range("M89").Value = -3749
$ws.Range("N89").Value = -22454.7785
$ws.Range("H112").Value = 1623.5625
$ws.Range("J112").Value = 2533.3333
$ws.Range("L112").Value = 7599.999899999999
$ws.Range("N112").Value = -9815.999899999999
$ws.Range("H118").Value = 165.25
$ws.Range("I118").Value = 165.25
$ws.Range("K118").Value = 495.75
$ws.Range("M118").Value = 1161.25
$ws.Range("H138").Value = 2420.1667
$ws.Range("J138").Value = 3270.8235
$ws.Range("L138").Value = 9812.470499999999
$ws.Range("N138").Value = -20092.4705

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4197.1763
$ws.Range("I32").Value = 3446.4167
$ws.Range("K32").Value = 3446.4167
$ws.Range("M32").Value = -3159.4167
$ws.Range("H61").Value = 9259.700000000001
$ws.Range("I61").Value = 6299.6665
$ws.Range("K61").Value = 6299.6665
$ws.Range("M61").Value = -6087.6665
$ws.Range("H63").Value = 2698.8
$ws.Range("J63").Value = 2447
$ws.Range("L63").Value = 2447
$ws.Range("N63").Value = -3819
$ws.Range("H66").Value = 2698.8
$ws.Range("J66").Value = 2447
$ws.Range("L66").Value = 12235
$ws.Range("N66").Value = -19099
$ws.Range("H74").Value = 25643930
$ws.Range("I74").Value = 37040344
$ws.Range("K74").Value = 37040344
$ws.Range("M74").Value = -37039470
$ws.Range("H77").Value = 25643930
$ws.Range("I77").Value = 37040344
$ws.Range("K77").Value = 185201720
$ws.Range("M77").Value = -185197352
$ws.Range("H114").Value = 69999
$ws.Range("J114").Value = 69999
$ws.Range("L114").Value = 69999
$ws.Range("N114").Value = -78677
$ws.Range("H136").Value = 9259.700000000001
$ws.Range("I136").Value = 6299.6665
$ws.Range("K136").Value = 18898.9995
$ws.Range("M136").Value = -16348.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5340.357
$ws.Range("I86").Value = 2975.9
$ws.Range("K86").Value = 2975.9
$ws.Range("M86").Value = -1852.9
$ws.Range("H89").Value = 5340.357
$ws.Range("I89").Value = 2975.9
$ws.Range("K89").Value = 14879.5
$ws.Range("M89").Value = -9263.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 60436.473
$ws.Range("I31").Value = 7752.7144
$ws.Range("J31").Value = 91168.664
$ws.Range("K31").Value = 7752.7144
$ws.Range("L31").Value = 91168.664
$ws.Range("M31").Value = -7457.7144
$ws.Range("N31").Value = -91758.664
$ws.Range("H34").Value = 60436.473
$ws.Range("I34").Value = 7752.7144
$ws.Range("J34").Value = 91168.664
$ws.Range("K34").Value = 7752.7144
$ws.Range("L34").Value = 91168.664
$ws.Range("M34").Value = -7550.7144
$ws.Range("N34").Value = -91572.664
$ws.Range("H107").Value = 971.381
$ws.Range("I107").Value = 956.1
$ws.Range("K107").Value = 956.1
$ws.Range("M107").Value = 963.9
$ws.Range("H122").Value = 5422.7856
$ws.Range("I122").Value = 1852.0714
$ws.Range("J122").Value = 8993.5
$ws.Range("K122").Value = 5556.2142
$ws.Range("L122").Value = 26980.5
$ws.Range("M122").Value = -3106.2142
$ws.Range("N122").Value = -31880.5
$ws.Range("H141").Value = 173032.38
$ws.Range("J141").Value = 173032.38
$ws.Range("L141").Value = 173032.38
$ws.Range("N141").Value = -183392.38

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8456003
$ws.Range("J131").Value = 13890510
$ws.Range("L131").Value = 41671530
$ws.Range("N131").Value = -41681610

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H109").Value = 74356.336
$ws.Range("I109").Value = 73500
$ws.Range("J109").Value = 74784.5
$ws.Range("K109").Value = 73500
$ws.Range("L109").Value = 74784.5
$ws.Range("M109").Value = -72460
$ws.Range("N109").Value = -76864.5
$ws.Range("H117").Value = 33750
$ws.Range("J117").Value = 33750
$ws.Range("L117").Value = 33750
$ws.Range("N117").Value = -40634
$ws.Range("H122").Value = 3408.2632
$ws.Range("I122").Value = 2691.1177
$ws.Range("J122").Value = 9504
$ws.Range("K122").Value = 8073.353099999999
$ws.Range("L122").Value = 28512
$ws.Range("M122").Value = -5623.353099999999
$ws.Range("N122").Value = -33412

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6571.7144
$ws.Range("I22").Value = 3000
$ws.Range("J22").Value = 11334
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 11334
$ws.Range("M22").Value = -2705
$ws.Range("N22").Value = -11924
$ws.Range("H27").Value = 6571.7144
$ws.Range("I27").Value = 3000
$ws.Range("J27").Value = 11334
$ws.Range("K27").Value = 3000
$ws.Range("L27").Value = 11334
$ws.Range("M27").Value = -2893
$ws.Range("N27").Value = -11548
$ws.Range("H40").Value = 7481.52
$ws.Range("I40").Value = 7343.579
$ws.Range("J40").Value = 7918.3335
$ws.Range("K40").Value = 7343.579
$ws.Range("L40").Value = 7918.3335
$ws.Range("M40").Value = -7207.579
$ws.Range("N40").Value = -8190.3335
$ws.Range("H42").Value = 39013.25
$ws.Range("I42").Value = 39012.5
$ws.Range("J42").Value = 39014
$ws.Range("K42").Value = 39012.5
$ws.Range("L42").Value = 39014
$ws.Range("M42").Value = -38449.5
$ws.Range("N42").Value = -40140
$ws.Range("H49").Value = 39013.25
$ws.Range("I49").Value = 39012.5
$ws.Range("J49").Value = 39014
$ws.Range("K49").Value = 39012.5
$ws.Range("L49").Value = 39014
$ws.Range("M49").Value = -38865.5
$ws.Range("N49").Value = -39308
$ws.Range("H61").Value = 6325.5
$ws.Range("I61").Value = 1400.8
$ws.Range("J61").Value = 9843.143
$ws.Range("K61").Value = 1400.8
$ws.Range("L61").Value = 9843.143
$ws.Range("M61").Value = -1198.8
$ws.Range("N61").Value = -10247.143
$ws.Range("H113").Value = 6325.5
$ws.Range("I113").Value = 1400.8
$ws.Range("J113").Value = 9843.143
$ws.Range("K113").Value = 1400.8
$ws.Range("L113").Value = 9843.143
$ws.Range("M113").Value = 769.2
$ws.Range("N113").Value = -14183.143
$ws.Range("H122").Value = 9080
$ws.Range("I122").Value = 7868.75
$ws.Range("K122").Value = 23606.25
$ws.Range("M122").Value = -21156.25
$ws.Range("H132").Value = 2914.4883
$ws.Range("J132").Value = 5536.615
$ws.Range("L132").Value = 16609.845
$ws.Range("N132").Value = -21669.845

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30298
$ws.Range("H96").Value = 4605.75
$ws.Range("I96").Value = 3538
$ws.Range("K96").Value = 3538
$ws.Range("M96").Value = -2165
$ws.Range("H100").Value = 1130
$ws.Range("J100").Value = 3003
$ws.Range("L100").Value = 6006
$ws.Range("N100").Value = -7088
$ws.Range("H119").Value = 78943.5
$ws.Range("J119").Value = 78943.5
$ws.Range("L119").Value = 78943.5
$ws.Range("N119").Value = -88619.5
$ws.Range("H136").Value = 3569.25
$ws.Range("I136").Value = 2664.5
$ws.Range("J136").Value = 9902.5
$ws.Range("K136").Value = 7993.5
$ws.Range("L136").Value = 29707.5
$ws.Range("M136").Value = -5443.5
$ws.Range("N136").Value = -34807.5

Write-Host "Applied $([string]259) cell updates across 8 sheets."
